$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates (rich-text concatenated as plain text; ---
# --- all runs share identical formatting, so this is visually identical) ---
$ws.Range("A8").Value = "Volume 29   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/21/2022  Through  11/27/2022"

# --- Style+type transplants (text <-> number) ---
# Donor cells (never modified by this script, used only to "borrow" a cell style):
#   C14 -> style s="14", shared text v="20" ("0")
#   E14 -> style s="14", shared text v="21" ("***.*")
#   F14 -> style s="15" (numeric #,##0)
#   K14 -> style s="16" (numeric #,##0.0)
$ws.Range("C14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))
$ws.Range("F14").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("K14").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("F14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 2
$ws.Range("F14").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 3
$ws.Range("K14").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100
$ws.Range("F14").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1
$ws.Range("F14").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("K14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = 0
$ws.Range("F14").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 1
$ws.Range("F14").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("K14").Copy($ws.Range("E29"))
$ws.Range("E29").Value = 0

# --- Plain value updates (style/type unchanged) ---
$ws.Range("L14").Value = 15.384615384615
$ws.Range("J15").Value = 26
$ws.Range("K15").Value = 7.692307692307
$ws.Range("M15").Value = 27.272727272727
$ws.Range("F16").Value = 37
$ws.Range("G16").Value = 47
$ws.Range("H16").Value = -21.276595744680
$ws.Range("I16").Value = 505
$ws.Range("J16").Value = 388
$ws.Range("K16").Value = 30.154639175257
$ws.Range("L16").Value = 23.774509803921
$ws.Range("M16").Value = 22.572815533980
$ws.Range("N16").Value = -71.159337521416
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = 7.692307692307
$ws.Range("F17").Value = 60
$ws.Range("G17").Value = 63
$ws.Range("H17").Value = -4.761904761904
$ws.Range("I17").Value = 767
$ws.Range("J17").Value = 647
$ws.Range("K17").Value = 18.547140649149
$ws.Range("L17").Value = 15.512048192771
$ws.Range("M17").Value = 92.713567839196
$ws.Range("N17").Value = -25.09765625
$ws.Range("C18").Value = 8
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 32
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = 60
$ws.Range("I18").Value = 340
$ws.Range("J18").Value = 170
$ws.Range("K18").Value = 100
$ws.Range("L18").Value = 31.782945736434
$ws.Range("M18").Value = 83.783783783783
$ws.Range("N18").Value = -72.380178716490
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 14
$ws.Range("F19").Value = 48
$ws.Range("H19").Value = -25
$ws.Range("I19").Value = 674
$ws.Range("J19").Value = 668
$ws.Range("K19").Value = 0.898203592814
$ws.Range("L19").Value = 1.812688821752
$ws.Range("M19").Value = 71.501272264631
$ws.Range("N19").Value = -4.125177809388
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 66.666666666666
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 90
$ws.Range("I20").Value = 281
$ws.Range("J20").Value = 164
$ws.Range("K20").Value = 71.341463414634
$ws.Range("L20").Value = 73.456790123456
$ws.Range("M20").Value = 153.153153153153
$ws.Range("N20").Value = -53.782894736842
$ws.Range("C21").Value = 48
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = 9.090909090909
$ws.Range("F21").Value = 199
$ws.Range("G21").Value = 206
$ws.Range("H21").Value = -3.398058252427
$ws.Range("I21").Value = 2610
$ws.Range("J21").Value = 2078
$ws.Range("K21").Value = 25.601539942252
$ws.Range("L21").Value = 18.906605922551
$ws.Range("M21").Value = 70.032573289902
$ws.Range("N21").Value = -52.171522814733
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 5
$ws.Range("H22").Value = -37.5
$ws.Range("I22").Value = 80
$ws.Range("J22").Value = 65
$ws.Range("K22").Value = 23.076923076923
$ws.Range("L22").Value = 12.676056338028
$ws.Range("M22").Value = 42.857142857142
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -11.111111111111
$ws.Range("F23").Value = 33
$ws.Range("H23").Value = 32
$ws.Range("I23").Value = 359
$ws.Range("J23").Value = 268
$ws.Range("K23").Value = 33.955223880597
$ws.Range("L23").Value = 25.087108013937
$ws.Range("M23").Value = 29.602888086642
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = -29.268292682926
$ws.Range("F24").Value = 128
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = -3.030303030303
$ws.Range("I24").Value = 1741
$ws.Range("J24").Value = 1252
$ws.Range("K24").Value = 39.057507987220
$ws.Range("L24").Value = 19.574175824175
$ws.Range("M24").Value = 35.486381322957
$ws.Range("C25").Value = 14
$ws.Range("E25").Value = -36.363636363636
$ws.Range("F25").Value = 68
$ws.Range("H25").Value = -24.444444444444
$ws.Range("I25").Value = 941
$ws.Range("J25").Value = 850
$ws.Range("K25").Value = 10.705882352941
$ws.Range("L25").Value = -1.362683438155
$ws.Range("M25").Value = -11.310084825636
$ws.Range("J26").Value = 65
$ws.Range("K26").Value = -29.230769230769
$ws.Range("L26").Value = -2.127659574468
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 78
$ws.Range("J27").Value = 83
$ws.Range("K27").Value = -6.024096385542
$ws.Range("L27").Value = 1.298701298701
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 66.666666666666
$ws.Range("I28").Value = 57
$ws.Range("J28").Value = 65
$ws.Range("K28").Value = -12.307692307692
$ws.Range("L28").Value = -13.636363636363
$ws.Range("M28").Value = 1.785714285714
$ws.Range("N28").Value = -71.921182266009
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 33.333333333333
$ws.Range("I29").Value = 47
$ws.Range("J29").Value = 57
$ws.Range("K29").Value = -17.543859649122
$ws.Range("L29").Value = -11.320754716981
$ws.Range("M29").Value = -2.083333333333
$ws.Range("N29").Value = -74.316939890710

Write-Host "Applied weekly crime data update"
